$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7009156306172883
$ws.Range("C2").Value = 0.2860586664313303
$ws.Range("D2").Value = 0.07717596355745115
$ws.Range("E2").Value = 0.09221461245308404
$ws.Range("G2").Value = 0.002539973165283912
$ws.Range("I2").Value = 1.600881764422098
$ws.Range("K2").Value = 0.984777650065439
$ws.Range("M2").Value = 0.3624649866762653
$ws.Range("N2").Value = 3.124431844195726

$ws.Range("B3").Value = 0.6577061489604148
$ws.Range("C3").Value = 0.2687370489676368
$ws.Range("D3").Value = 0.07019141311180022
$ws.Range("E3").Value = 0.0837454226861496
$ws.Range("G3").Value = 0.002544871371795657
$ws.Range("I3").Value = 1.578277073721466
$ws.Range("K3").Value = 0.9244056624362713
$ws.Range("M3").Value = 0.3362313141047224
$ws.Range("N3").Value = 3.105790780679456

$ws.Range("B4").Value = 0.6316434915227092
$ws.Range("C4").Value = 0.2582767995685629
$ws.Range("D4").Value = 0.06594376253887901
$ws.Range("E4").Value = 0.07859841047611837
$ws.Range("G4").Value = 0.00254803530986642
$ws.Range("I4").Value = 1.564986640784312
$ws.Range("K4").Value = 0.8879779473368785
$ws.Range("M4").Value = 0.3203529005920416
$ws.Range("N4").Value = 3.094969074603398

$ws.Range("B5").Value = 0.6211397965495564
$ws.Range("C5").Value = 0.2540578089431449
$ws.Range("D5").Value = 0.06422295955457002
$ws.Range("E5").Value = 0.07651403418850578
$ws.Range("G5").Value = 0.002549364109657004
$ws.Range("I5").Value = 1.559718239165534
$ws.Range("K5").Value = 0.8732934115564035
$ws.Range("M5").Value = 0.3139393013950738
$ws.Range("N5").Value = 3.090715380239956

$ws.Range("B6").Value = 0.6194027162623286
$ws.Range("C6").Value = 0.2533598723234434
$ws.Range("D6").Value = 0.06393783064599745
$ws.Range("E6").Value = 0.0761687065186436
$ws.Range("G6").Value = 0.002549587143475626
$ws.Range("I6").Value = 1.558852321727514
$ws.Range("K6").Value = 0.8708646902901762
$ws.Range("M6").Value = 0.312877751152449
$ws.Range("N6").Value = 3.090018474604051

$ws.Range("B7").Value = 0.6315013619672527
$ws.Range("C7").Value = 0.2582197246337614
$ws.Range("D7").Value = 0.06592051426878243
$ws.Range("E7").Value = 0.07857024728351547
$ws.Range("G7").Value = 0.002548053070362483
$ws.Range("I7").Value = 1.564914992543123
$ws.Range("K7").Value = 0.8877792599599275
$ws.Range("M7").Value = 0.3202661746693281
$ws.Range("N7").Value = 3.094911076103372

$ws.Range("B8").Value = 0.6859194189591165
$ws.Range("C8").Value = 0.2800495226738349
$ws.Range("D8").Value = 0.07475912063443957
$ws.Range("E8").Value = 0.08928324187217385
$ws.Range("G8").Value = 0.002541629684612313
$ws.Range("I8").Value = 1.592965098817416
$ws.Range("K8").Value = 0.9638276687254859
$ws.Range("M8").Value = 0.3533716552823165
$ws.Range("N8").Value = 3.117874486598623

$ws.Range("B9").Value = 0.7963843966533375
$ws.Range("C9").Value = 0.3242706862374405
$ws.Range("D9").Value = 0.09242303307460986
$ws.Range("E9").Value = 0.1107265402809574
$ws.Range("G9").Value = 0.002530268321506234
$ws.Range("I9").Value = 1.652672748607358
$ws.Range("K9").Value = 1.118103333803958
$ws.Range("M9").Value = 0.4201425368342413
$ws.Range("N9").Value = 3.167890545111902

$ws.Range("B10").Value = 0.879890695803681
$ws.Range("C10").Value = 0.3576557301757362
$ws.Range("D10").Value = 0.105613798190646
$ws.Range("E10").Value = 0.1267669343037809
$ws.Range("G10").Value = 0.002522665185442518
$ws.Range("I10").Value = 1.699451874958839
$ws.Range("K10").Value = 1.234681923594195
$ws.Range("M10").Value = 0.470378495834197
$ws.Range("N10").Value = 3.207728619333182

$ws.Range("B11").Value = 0.9184035392489136
$ws.Range("C11").Value = 0.3730456772863988
$ws.Range("D11").Value = 0.1116633739582795
$ws.Range("E11").Value = 0.13413083683173
$ws.Range("G11").Value = 0.002519366015998145
$ws.Range("I11").Value = 1.721375604060086
$ws.Range("K11").Value = 1.288440135546523
$ws.Range("M11").Value = 0.4934996555785318
$ws.Range("N11").Value = 3.226534983860091

$ws.Range("B12").Value = 0.9330638539443612
$ws.Range("C12").Value = 0.378903190686998
$ws.Range("D12").Value = 0.1139614163458162
$ws.Range("E12").Value = 0.1369293418056543
$ws.Range("G12").Value = 0.002518139504348323
$ws.Range("I12").Value = 1.729770844907065
$ws.Range("K12").Value = 1.308902895100573
$ws.Range("M12").Value = 0.5022944860015457
$ws.Range("N12").Value = 3.233755659162938

$ws.Range("B13").Value = 0.9299030851225893
$ws.Range("C13").Value = 0.3776403426398076
$ws.Range("D13").Value = 0.113466169217233
$ws.Range("E13").Value = 0.1363261859035489
$ws.Range("G13").Value = 0.002518402643169456
$ws.Range("I13").Value = 1.727958624122351
$ws.Range("K13").Value = 1.304491152262983
$ws.Range("M13").Value = 0.5003986001267009
$ws.Range("N13").Value = 3.232196136909494

$ws.Range("B14").Value = 0.9196081177144038
$ws.Range("C14").Value = 0.3735269808596513
$ws.Range("D14").Value = 0.1118522901422381
$ws.Range("E14").Value = 0.1343608701976464
$ws.Range("G14").Value = 0.002519264653975244
$ws.Range("I14").Value = 1.722064413986786
$ws.Range("K14").Value = 1.290121493405366
$ws.Range("M14").Value = 0.4942224187420123
$ws.Range("N14").Value = 3.22712704228249

$ws.Range("B15").Value = 0.9133121177411567
$ws.Range("C15").Value = 0.3710113072905585
$ws.Range("D15").Value = 0.1108646846279839
$ws.Range("E15").Value = 0.1331583633821865
$ws.Range("G15").Value = 0.002519795627338617
$ws.Range("I15").Value = 1.718466196769668
$ws.Range("K15").Value = 1.281333465260616
$ws.Range("M15").Value = 0.490444475540329
$ws.Range("N15").Value = 3.224035006250148

$ws.Range("B16").Value = 0.8773844300973224
$ws.Range("C16").Value = 0.3566540863519947
$ws.Range("D16").Value = 0.1052194434788447
$ws.Range("E16").Value = 0.1262870618283429
$ws.Range("G16").Value = 0.002522883993418558
$ws.Range("I16").Value = 1.698032114559595
$ws.Range("K16").Value = 1.23118341723216
$ws.Range("M16").Value = 0.4688729400487972
$ws.Range("N16").Value = 3.206513405919651

$ws.Range("B17").Value = 0.8554790797299177
$ws.Range("C17").Value = 0.3478987191203657
$ws.Range("D17").Value = 0.1017689340091863
$ws.Range("E17").Value = 0.1220891350813176
$ws.Range("G17").Value = 0.00252481937491034
$ws.Range("I17").Value = 1.68566177907735
$ws.Range("K17").Value = 1.200604850138205
$ws.Range("M17").Value = 0.4557088041667043
$ws.Range("N17").Value = 3.1959401671877

$ws.Range("B18").Value = 0.8429291077411278
$ws.Range("C18").Value = 0.3428819335242963
$ws.Range("D18").Value = 0.09978889819039694
$ws.Range("E18").Value = 0.1196808904132709
$ws.Range("G18").Value = 0.002525947579049004
$ws.Range("I18").Value = 1.678607235507215
$ws.Range("K18").Value = 1.183085125647153
$ws.Range("M18").Value = 0.4481623762129985
$ws.Range("N18").Value = 3.189923037478707

$ws.Range("B19").Value = 0.8386883648070409
$ws.Range("C19").Value = 0.3411865964874039
$ws.Range("D19").Value = 0.09911927866292558
$ws.Range("E19").Value = 0.1188665712752979
$ws.Range("G19").Value = 0.002526332153972354
$ws.Range("I19").Value = 1.676229067770066
$ws.Range("K19").Value = 1.177164932522061
$ws.Range("M19").Value = 0.4456115974481065
$ws.Range("N19").Value = 3.187896763738621

$ws.Range("B20").Value = 0.857805822698964
$ws.Range("C20").Value = 0.3488287656707882
$ws.Range("D20").Value = 0.1021357688579485
$ws.Range("E20").Value = 0.1225353581117403
$ws.Range("G20").Value = 0.002524611796130173
$ws.Range("I20").Value = 1.686972351524531
$ws.Range("K20").Value = 1.203852917040223
$ws.Range("M20").Value = 0.4571075319427536
$ws.Range("N20").Value = 3.197059044320838

$ws.Range("B21").Value = 0.9226299211144351
$ws.Range("C21").Value = 0.3747343655897737
$ws.Range("D21").Value = 0.1123261289641562
$ws.Range("E21").Value = 0.1349378583673513
$ws.Range("G21").Value = 0.002519010842476419
$ws.Range("I21").Value = 1.723793151777187
$ws.Range("K21").Value = 1.294339331707477
$ws.Range("M21").Value = 0.4960354391418491
$ws.Range("N21").Value = 3.228613262075072

$ws.Range("B22").Value = 0.9654415408225532
$ws.Range("C22").Value = 0.3918383730051858
$ws.Range("D22").Value = 0.119028152188406
$ws.Range("E22").Value = 0.1431017825500618
$ws.Range("G22").Value = 0.00251548320319064
$ws.Range("I22").Value = 1.748401222804247
$ws.Range("K22").Value = 1.354094100403017
$ws.Range("M22").Value = 0.5217067793346075
$ws.Range("N22").Value = 3.249813844807818

$ws.Range("B23").Value = 0.9425511633793349
$ws.Range("C23").Value = 0.3826936291443133
$ws.Range("D23").Value = 0.1154472604524983
$ws.Range("E23").Value = 0.1387391150039363
$ws.Range("G23").Value = 0.002517353851850067
$ws.Range("I23").Value = 1.735217482362003
$ws.Range("K23").Value = 1.322144998475608
$ws.Range("M23").Value = 0.5079842403816883
$ws.Range("N23").Value = 3.238445540133768

$ws.Range("B24").Value = 0.8567537661267011
$ws.Range("C24").Value = 0.3484082393619303
$ws.Range("D24").Value = 0.1019699112873553
$ws.Range("E24").Value = 0.1223336045114394
$ws.Range("G24").Value = 0.002524705594185311
$ws.Range("I24").Value = 1.686379663192412
$ws.Range("K24").Value = 1.20238427826763
$ws.Range("M24").Value = 0.4564750992436331
$ws.Range("N24").Value = 3.196553008231803

$ws.Range("B25").Value = 0.766092861513755
$ws.Range("C25").Value = 0.3121528939929021
$ws.Range("D25").Value = 0.08760786066839898
$ws.Range("E25").Value = 0.1048768455206357
$ws.Range("G25").Value = 0.002533210576463487
$ws.Range("I25").Value = 1.636012230179986
$ws.Range("K25").Value = 1.075807058474965
$ws.Range("M25").Value = 0.4018758443245645
$ws.Range("N25").Value = 3.153820834207124
